# Scheduled-runner refresh of computed Leve profit columns (H-N) across all 8 job sheets.
# Values below come from an external price-fetch recompute; row layout is unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 550.46155
$ws.Range("J12").Value = 758.6667
$ws.Range("L12").Value = 758.6667
$ws.Range("N12").Value = -1098.6667

# Row 75: Tomes Roam on the Range / Dhalmelskin Codex
$ws.Range("H75").Value = 67653
$ws.Range("J75").Value = 67653
$ws.Range("L75").Value = 67653
$ws.Range("N75").Value = -69525

# Row 78: Field Trip to the Unknown (L) / Dhalmelskin Codex
$ws.Range("H78").Value = 67653
$ws.Range("J78").Value = 67653
$ws.Range("L78").Value = 202959
$ws.Range("N78").Value = -212319

# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 7359.2856
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 7359.2856
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 22077.8568
$ws.Range("N80").Value = -24073.8568
$ws.Range("M80").ClearContents()

# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 7359.2856
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 7359.2856
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 66233.5704
$ws.Range("N83").Value = -76217.5704
$ws.Range("M83").ClearContents()

# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 639
$ws.Range("I107").Value = 645.5833
$ws.Range("J107").Value = 599.5
$ws.Range("K107").Value = 645.5833
$ws.Range("L107").Value = 599.5
$ws.Range("M107").Value = 1274.4167
$ws.Range("N107").Value = -4439.5

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 57185.668
$ws.Range("I132").Value = 57185.668
$ws.Range("K132").Value = 171557.004
$ws.Range("M132").Value = -169027.004

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 4904.1953
$ws.Range("I2").Value = 5384.8335
$ws.Range("K2").Value = 5384.8335
$ws.Range("M2").Value = -5271.8335

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 4904.1953
$ws.Range("I116").Value = 5384.8335
$ws.Range("K116").Value = 5384.8335
$ws.Range("M116").Value = -3090.8335

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2007.3103
$ws.Range("I122").Value = 2055.3809
$ws.Range("K122").Value = 6166.1427
$ws.Range("M122").Value = -3716.1427

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 20005180
$ws.Range("J132").Value = 50008756
$ws.Range("L132").Value = 150026268
$ws.Range("N132").Value = -150031328

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 4904.1953
$ws.Range("I3").Value = 5384.8335
$ws.Range("K3").Value = 5384.8335
$ws.Range("M3").Value = -5270.8335

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2666.6667
$ws.Range("J86").Value = 2600
$ws.Range("L86").Value = 2600
$ws.Range("N86").Value = -4846

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2666.6667
$ws.Range("J89").Value = 2600
$ws.Range("L89").Value = 13000
$ws.Range("N89").Value = -24232

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 10872728
$ws.Range("I134").Value = 13892018
$ws.Range("K134").Value = 41676054
$ws.Range("M134").Value = -41673519

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1811.3334
$ws.Range("I31").Value = 1032.4166
$ws.Range("K31").Value = 1032.4166
$ws.Range("M31").Value = -737.4166

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1811.3334
$ws.Range("I34").Value = 1032.4166
$ws.Range("K34").Value = 1032.4166
$ws.Range("M34").Value = -830.4166

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1200
$ws.Range("I58").Value = 1200
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1200
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -997
$ws.Range("N58").ClearContents()

# Row 59: Bow Down to Magic / Crab Bow
$ws.Range("H59").Value = 31209.889
$ws.Range("J59").Value = 31209.889
$ws.Range("L59").Value = 31209.889
$ws.Range("N59").Value = -33499.889

# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 2950
$ws.Range("I86").Value = 2950
$ws.Range("K86").Value = 2950
$ws.Range("M86").Value = -1827

# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 2950
$ws.Range("I89").Value = 2950
$ws.Range("K89").Value = 14750
$ws.Range("M89").Value = -9134

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1200
$ws.Range("I136").Value = 1200
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3600
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1050
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 48: Rise and Dine / Cheese Souffle
$ws.Range("H48").Value = 10599.667
$ws.Range("J48").Value = 14499.5
$ws.Range("L48").Value = 43498.5
$ws.Range("N48").Value = -43998.5

# Row 69: Loving That Muffin Top / Ishgardian Muffin
$ws.Range("H69").Value = 3751.111
$ws.Range("I69").Value = 2400
$ws.Range("J69").Value = 4137.143
$ws.Range("K69").Value = 7200
$ws.Range("L69").Value = 12411.429
$ws.Range("M69").Value = -6389
$ws.Range("N69").Value = -14033.429

# Row 72: Muffin of the Morn (L) / Ishgardian Muffin
$ws.Range("H72").Value = 3751.111
$ws.Range("I72").Value = 2400
$ws.Range("J72").Value = 4137.143
$ws.Range("K72").Value = 21600
$ws.Range("L72").Value = 37234.287
$ws.Range("M72").Value = -17544
$ws.Range("N72").Value = -45346.287

$ws = $wb.Worksheets.Item("GSM")
# Row 63: Not on My Table / Mythrite Earrings of Healing
$ws.Range("H63").Value = 34971
$ws.Range("J63").Value = 34971
$ws.Range("L63").Value = 34971
$ws.Range("N63").Value = -36343

# Row 66: Heinz's Dilemma (L) / Mythrite Earrings of Healing
$ws.Range("H66").Value = 34971
$ws.Range("J66").Value = 34971
$ws.Range("L66").Value = 104913
$ws.Range("N66").Value = -111777

# Row 92: Play It by Ear / Triphane Earrings of Healing
$ws.Range("H92").Value = 1975
$ws.Range("J92").Value = 1975
$ws.Range("L92").Value = 1975
$ws.Range("N92").Value = -5719

# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 1613.0834
$ws.Range("I97").Value = 1434.25
$ws.Range("K97").Value = 1434.25
$ws.Range("M97").Value = -938.25

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 7241
$ws.Range("I126").Value = 10048.75
$ws.Range("K126").Value = 30146.25
$ws.Range("M126").Value = -27676.25

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 3112.5
$ws.Range("I132").Value = 3137
$ws.Range("K132").Value = 9411
$ws.Range("M132").Value = -6881

$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 1593.625
$ws.Range("I55").Value = 850
$ws.Range("J55").Value = 2833
$ws.Range("K55").Value = 850
$ws.Range("L55").Value = 2833
$ws.Range("M55").Value = -677
$ws.Range("N55").Value = -3179

# Row 58: Handle with Care / Peisteskin Cesti
$ws.Range("H58").Value = 18490.857
$ws.Range("I58").Value = 2334.5
$ws.Range("J58").Value = 40032.668
$ws.Range("K58").Value = 2334.5
$ws.Range("L58").Value = 40032.668
$ws.Range("M58").Value = -2074.5
$ws.Range("N58").Value = -40552.668

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 5381.357
$ws.Range("I82").Value = 3944.5
$ws.Range("K82").Value = 3944.5
$ws.Range("M82").Value = -3583.5

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 5381.357
$ws.Range("I85").Value = 3944.5
$ws.Range("K85").Value = 3944.5
$ws.Range("M85").Value = -2696.5

$ws = $wb.Worksheets.Item("WVR")
# Row 69: Fashion Patrol / Holy Rainbow Sarouel of Casting
$ws.Range("H69").Value = 48795
$ws.Range("J69").Value = 48795
$ws.Range("L69").Value = 48795
$ws.Range("N69").Value = -50293

# Row 72: Dress Code Violation (L) / Holy Rainbow Sarouel of Casting
$ws.Range("H72").Value = 48795
$ws.Range("J72").Value = 48795
$ws.Range("L72").Value = 146385
$ws.Range("N72").Value = -153873

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 3042
$ws.Range("I126").Value = 2458
$ws.Range("J126").Value = 3366.4443
$ws.Range("K126").Value = 7374
$ws.Range("L126").Value = 10099.3329
$ws.Range("M126").Value = -4904
$ws.Range("N126").Value = -15039.3329

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 704.2
$ws.Range("I136").Value = 704.2
$ws.Range("K136").Value = 2112.6
$ws.Range("M136").Value = 437.3999999999996
